$wb = $excel.ActiveWorkbook

$wsAll = $wb.Worksheets.Item("allFields")
$wsMissing = $wb.Worksheets.Item("missingFields")

# --- Fix the two multi-line error messages in "missingFields" so the
# --- literal "\n" text is replaced by a real line break, then make the
# --- cells wrap so the break is visible, and grow the rows to fit.
$wsMissing.Range("H6").Value = "Lastname can't be blank`nUserName can't be blank"
$wsMissing.Range("H6").WrapText = $true
$wsMissing.Rows.Item(6).RowHeight = 87

$wsMissing.Range("H7").Value = "Firstname can't be blank`nLastname can't be blank`nUserName can't be blank`nPasswd can't be blank"
$wsMissing.Range("H7").WrapText = $true
$wsMissing.Rows.Item(7).RowHeight = 174

# --- Move the active tab from "allFields" to "missingFields" and restore
# --- the selection on each sheet.
$wsAll.Range("G3").Select()

$wsMissing.Activate()
$wsMissing.Range("H6:H7").Select()
